# Applies the "Updated symbol list" commit: refreshed prices for several
# coins, plus several rows shifted up one coin (the row keeps its rank
# label in column E but the coin/link/price in B/C/D move to the row
# above), and a couple of "Worstin24h"/"Bestin24h" suffix tweaks in E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of column D in this sheet is stored as TEXT (e.g. "249.15"), not a
# number, even though it looks numeric. Plain `.Value = "250.07"` would get
# auto-coerced to a real number by Excel, so force text storage first via
# NumberFormat "@", then strip the format back off so we don't leave a
# stray text-format style behind on cells that had none originally.
$priceCells = @(
  "D2","D3","D4","D5","D7","D8","D9",
  "D10","D11","D12","D13","D14","D15","D16","D17","D18",
  "D19","D20","D21","D22","D23","D24","D25","D26",
  "D40","D41","D42","D43","D44","D45","D48","D49"
)
foreach ($addr in $priceCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value  = "250.07"
$ws.Range("D3").Value  = "22.69"
$ws.Range("D4").Value  = "5.381"
$ws.Range("D5").Value  = "0.05620"
$ws.Range("D7").Value  = "6.360"
$ws.Range("D8").Value  = "0.8152"
$ws.Range("D9").Value  = "0.9190"

$ws.Range("D19").Value = "0.006396"
$ws.Range("D20").Value = "0.004997"
$ws.Range("D21").Value = "0.001030"
$ws.Range("D22").Value = "0.0001499"
$ws.Range("D23").Value = "3.725"
$ws.Range("D24").Value = "2.161"
$ws.Range("D25").Value = "0.3252"
$ws.Range("D26").Value = "0.1303"

$ws.Range("D40").Value = "0.03973"
$ws.Range("D44").Value = "0.007538"
$ws.Range("D45").Value = "0.00005572"
$ws.Range("D48").Value = "0.6752"
$ws.Range("D49").Value = "0.2220"

# Rows 10-18: each row's Coin/Link/Price take on the identity of the row
# that used to be below it (One -> WazirX -> MandalaExchangeToken -> ... ->
# CoinExToken -> One), while E keeps its rank number and gets the new name.
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1422"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.07487"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "0.03191"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03089"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09328"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "3.557"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001633"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04715"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0005760"
$ws.Range("E18").Value = "17OneONE"

# Row 27: AAXToken keeps its coin/link/price, only the "Worstin24h" flag in
# E is dropped.
$ws.Range("E27").Value = "26AAXTokenAAB"

# Rows 41-43: similar rotation (KickToken -> BKEXToken -> CEJI ->
# KickToken), and row 43 picks up a new "Worstin24h" flag in E.
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "0.1067"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.002728"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "0.002921"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

# Drop the temporary text NumberFormat so the D cells end up with no
# explicit style, matching their original (unstyled) state.
foreach ($addr in $priceCells) { $ws.Range($addr).ClearFormats() }

Write-Host "Applied symbol-list update to $($priceCells.Count) price cells and related rows"
